$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (shifting old Treatment column from B to C)
$ws.Columns("B").Insert()

# Fill new Stream column (B) with base names first, so new shared
# strings are appended in the same order as the original edit:
# CHUCK, LOON, MCTE, W-100, W-113, then RowName last.
$ws.Range("B2").Value = "CHUCK"
$ws.Range("B3").Value = "CHUCK"
$ws.Range("B4").Value = "LOON"
$ws.Range("B5").Value = "LOON"
$ws.Range("B6").Value = "MCTE"
$ws.Range("B7").Value = "MCTE"
$ws.Range("B8").Value = "W-100"
$ws.Range("B9").Value = "W-100"
$ws.Range("B10").Value = "W-113"
$ws.Range("B11").Value = "W-113"

# Header row
$ws.Range("B1").Value = "Stream"
$ws.Range("C1").Value = "Treatment"
$ws.Range("A1").Value = "RowName"
